$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new cell value for C4
$ws.Range("C4").Value = "ATTRIBUTE"

# Update the active selection to C4 (matches the new selection in the diff)
$ws.Range("C4").Select()
